# Add the ItemsBase header row to Лист1 (sheet 1): Key, Name, Type, Quality, Description
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header values ---
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Type"
$ws.Range("D1").Value = "Quality"
$ws.Range("E1").Value = "Description"

# --- A1 "Key": Check Cell style, font swapped for a custom display font ---
$ws.Range("A1").Style = "Check Cell"
$ws.Range("A1").Font.Name = "18thCentury"

# --- B1 "Name": Heading 1 style, font color turned red ---
$ws.Range("B1").Style = "Heading 1"
$ws.Range("B1").Font.Color = 255

# --- C1:E1 "Type"/"Quality"/"Description": Heading 1 style, font color turned purple ---
$ws.Range("C1").Style = "Heading 1"
$ws.Range("C1").Font.Color = 10498160

$ws.Range("D1").Style = "Heading 1"
$ws.Range("D1").Font.Color = 10498160

$ws.Range("E1").Style = "Heading 1"
$ws.Range("E1").Font.Color = 10498160

# --- Give row 2 the same column span as row 1 (keeps the thick/double borders visually closed) ---
$ws.Rows.Item(2).Style = "Normal"

# --- Row heights matching the bigger header font + thick border ---
$ws.Rows.Item(1).RowHeight = 21
$ws.Rows.Item(2).RowHeight = 15.75

# --- Page setup: portrait A4/letter-ish "9" = A4 paper, portrait orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ends on E1, matching the saved cursor position ---
$ws.Range("E1").Select() | Out-Null
